$wb = $excel.ActiveWorkbook

# Update the "Ready for handoff" status text to "In Translation" everywhere it appears.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# Narrow the Status columns (was width 17.2159881591797 -> now 13.4101848602295 in the
# saved OOXML). Excel's COM ColumnWidth is quantized to whole pixels using the
# workbook's Normal-style max-digit-width (6px here) plus 5px of padding, i.e.
# storedWidth = (round(cw*MDW) + 5) / MDW. cw = 12.5 is the value that rounds to the
# closest reachable stored width (13.333333333333334) to the target.
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5

$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
